# fix(offer-problem): added an additional check when checking the price,
# if there are a small number of offers on the product.
#
# Appends a new product row (name + URL) to the "Products" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the URL (column B) before the product name (column A) so the two
# new shared-string entries land in the same order as the source workbook
# (URL first, then name).
$ws.Range("B6").Value = "https://www.emag.ro/prelata-acoperire-piscina-pvc-neagra-366-cm-bestway-8050060/pd/D828Z4MBM"
$ws.Range("A6").Value = "Prelata acoperire piscina, PVC, neagra, 366 cm, Bestway"

# Reflect the new selection state (user clicked the new row, then moved on
# to the next empty one).
$ws.Range("A6").Select()
$ws.Range("A15").Select()
